$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# New header cells on row 7 of the first sheet, typed in the same order
# the shared-string table grew in the target file (L7, M7, J7, K7).
$ws1.Range("L7").Value = "vr/v0"
$ws1.Range("M7").Value = "d(vr/v0)"
$ws1.Range("J7").Value = "d1(dvr/v0)"
$ws1.Range("K7").Value = "d2(vr*dv0/v0^2)"

# New calculated row underneath.
$ws1.Range("J8").Formula = "=C8/D8"
$ws1.Range("K8").Formula = "=B8*E8/(D8^2)"
$ws1.Range("L8").Formula = "=B8/D8"
$ws1.Range("M8").Formula = "=SQRT(SUM(J8^2,K8^2))"

# Widen the new/affected columns so the plotted values are readable
# (closest attainable widths to the target 12.625 / 13.375 "best fit"
# character widths given this host's column-width quantisation).
$ws1.Columns.Item(10).ColumnWidth = 11.8
$ws1.Columns.Item(11).ColumnWidth = 12.7

# Restore the selection on the second sheet first (so it is no longer the
# active tab once the first sheet is selected below).
$ws2.Range("H13").Select()

# Make the first sheet the active tab with O8 selected.
$ws1.Range("O8").Select()
